# Hortaliza, Feria Lagunitas de Puerto Montt - Brócoli
# Weekly update: insert two new price records (rows 160 and 252 in the
# final layout), pushing all subsequent rows down.
#
# Row 160 (new) sits right before the old row 160 (now row 161).
# Row 252 (new) sits right before what is, after the first insert, row 252
# (the old row 251, now shifted to row 253).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two blank rows -------------------------------------------
$ws.Rows("160:160").Insert()
$ws.Rows("252:252").Insert()

# --- Fill in the first new record (row 160) -------------------------------
$ws.Range("A160").Value = 4
$ws.Range("B160").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C160").Value = "Los Lagos"
$ws.Range("D160").Value = 44567
$ws.Range("E160").Value = 10
$ws.Range("F160").Value = 100112023
$ws.Range("G160").Value = "Brócoli"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 500
$ws.Range("K160").Value = 1400
$ws.Range("L160").Value = 1400
$ws.Range("M160").Value = 1400
$ws.Range("N160").Value = "$/unidad"
$ws.Range("O160").Value = "Región Metropolitana"
$ws.Range("P160").Value = 1400
$ws.Range("Q160").Value = 1
$ws.Range("R160").Value = "Hortaliza"

# --- Fill in the second new record (row 252) ------------------------------
$ws.Range("A252").Value = 4
$ws.Range("B252").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C252").Value = "Los Lagos"
$ws.Range("D252").Value = 44568
$ws.Range("E252").Value = 10
$ws.Range("F252").Value = 100112023
$ws.Range("G252").Value = "Brócoli"
$ws.Range("H252").Value = "Sin especificar"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 800
$ws.Range("K252").Value = 1300
$ws.Range("L252").Value = 1400
$ws.Range("M252").Value = 1350
$ws.Range("N252").Value = "$/unidad"
$ws.Range("O252").Value = "Región Metropolitana"
$ws.Range("P252").Value = 1350
$ws.Range("Q252").Value = 1
$ws.Range("R252").Value = "Hortaliza"
